# Generate Report for Handoff
# Replace the old localization file id (9f3e394b-f36c-4626-acb7-f1a4c235305e)
# with the new one (0c4274bb-7005-4af8-b3a2-b38817b0c895) across all sheets,
# including hyperlink display text, and refresh the associated timestamps
# and xliff file names.

$wb = $excel.ActiveWorkbook

$oldId = "9f3e394b-f36c-4626-acb7-f1a4c235305e"
$newId = "0c4274bb-7005-4af8-b3a2-b38817b0c895"

$oldHash = "3996be5b8af4d14f6f29a6647289901e08af3840"
$newHash = "085d18978db4ebc1024d9ed7103fcf0f50d9e06f"

# External hyperlink target is unchanged by this commit; reuse it when we
# recreate the hyperlinks below so the relationship target stays identical.
$hyperlinkAddress = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/132c797aea28bf5e2eab52d8a60bc1da075e044a/e2e/$oldId.md"

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = "$newId.md"

$rngB2 = $wsOverview.Range("B2")
$rngB2.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($rngB2, $hyperlinkAddress, "", "", "e2e\$newId.md")

$wsOverview.Range("G2").Value = "2016-09-01 07:09:16"

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$rngZhA2 = $wsZhCn.Range("A2")
$rngZhA2.Hyperlinks.Delete()
$wsZhCn.Hyperlinks.Add($rngZhA2, $hyperlinkAddress, "", "", "$newId.md")

$wsZhCn.Range("G2").Value = "$newId.$newHash.zh-cn.xlf"
$wsZhCn.Range("H2").Value = "2016-09-01 07:09:11"

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$rngDeA2 = $wsDeDe.Range("A2")
$rngDeA2.Hyperlinks.Delete()
$wsDeDe.Hyperlinks.Add($rngDeA2, $hyperlinkAddress, "", "", "$newId.md")

$wsDeDe.Range("G2").Value = "$newId.$newHash.de-de.xlf"
$wsDeDe.Range("H2").Value = "2016-09-01 07:09:16"
